# Update odds data and snapshot timestamps to reflect the latest Betfair
# pull for 2026-02-27 fixtures (commit: "Atualizando o arquivo XLSX").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTs = "2026-02-25 07:31:48"

# --- Row 5: Al-Hazm (KSA) vs Al-Ettifaq ---
$ws.Range("F5").Value = 2.82
$ws.Range("G5").Value = 3.3
$ws.Range("H5").Value = 2.32
$ws.Range("I5").Value = 2.82
$ws.Range("J5").Value = 3.4
$ws.Range("K5").Value = 4.1
$ws.Range("P5").Value = 2.02
$ws.Range("Q5").Value = 1.76
$ws.Range("BH5").Value = $newTs

# --- Row 6: Al-Ittihad vs Al-Khaleej Saihat ---
$ws.Range("G6").Value = 1.64
$ws.Range("H6").Value = 5.4
$ws.Range("K6").Value = 5.6
$ws.Range("BH6").Value = $newTs

# --- Row 7: Albacete vs Almeria ---
$ws.Range("P7").Value = 2.08
$ws.Range("Q7").Value = 1.75
$ws.Range("BH7").Value = $newTs

# --- Row 10: JS Saoura vs ES Ben Aknoun ---
$ws.Range("G10").Value = 1.8
$ws.Range("H10").Value = 2.3
$ws.Range("J10").Value = 2.24
$ws.Range("BH10").Value = $newTs

# --- Row 11: Cerro vs Boston River ---
$ws.Range("F11").Value = 3.2
$ws.Range("H11").Value = 2.48
$ws.Range("I11").Value = 3
$ws.Range("J11").Value = 2.66
$ws.Range("BH11").Value = $newTs

# --- Rows with only the snapshot timestamp refreshed ---
$ws.Range("BH2").Value = $newTs
$ws.Range("BH3").Value = $newTs
$ws.Range("BH4").Value = $newTs
$ws.Range("BH8").Value = $newTs
$ws.Range("BH9").Value = $newTs
